# Update Name of Algo
# Apply updated imputed values produced by the RandomForest algorithm
# to the corresponding data cells on the active worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = -7.461599999999999
$ws.Range("A3").Value = -22.22610000000001
$ws.Range("E3").Value = 15.7156
$ws.Range("E12").Value = 17.19270000000002
$ws.Range("A14").Value = -21.8824
$ws.Range("A21").Value = -19.97709999999998
$ws.Range("A23").Value = -20.40009999999998
$ws.Range("E24").Value = 16.5774
$ws.Range("A25").Value = -21.93290000000001
$ws.Range("D25").Value = -7.810399999999999
$ws.Range("E25").Value = 16.8265
$ws.Range("A26").Value = -21.05369999999997
$ws.Range("D27").Value = -8.581900000000005
$ws.Range("A29").Value = -20.94019999999998
$ws.Range("D31").Value = -8.610900000000008
$ws.Range("D39").Value = -7.828199999999998
$ws.Range("D48").Value = -7.204399999999997
$ws.Range("E50").Value = 16.2962
$ws.Range("D51").Value = -7.312999999999999
$ws.Range("D52").Value = -7.5933
$ws.Range("A53").Value = -21.9864
$ws.Range("E53").Value = 17.05690000000001
$ws.Range("D55").Value = -8.553499999999996
$ws.Range("D56").Value = -7.867200000000004
$ws.Range("A57").Value = -22.48540000000001
$ws.Range("D57").Value = -8.373299999999997
$ws.Range("E57").Value = 16.4861
$ws.Range("A59").Value = -22.23549999999999
$ws.Range("E61").Value = 16.55580000000001
$ws.Range("E63").Value = 17.45950000000002
$ws.Range("A69").Value = -21.6035
$ws.Range("E70").Value = 17.30100000000001
$ws.Range("D73").Value = -7.829899999999996
$ws.Range("A79").Value = -20.09860000000002
$ws.Range("A83").Value = -21.85340000000001
$ws.Range("E86").Value = 16.5707
$ws.Range("D89").Value = -5.733000000000003
$ws.Range("D90").Value = -8.080000000000002
$ws.Range("A91").Value = -21.41680000000002
$ws.Range("D92").Value = -5.845800000000001
$ws.Range("A93").Value = -21.01999999999999
$ws.Range("E98").Value = 15.84250000000001
$ws.Range("E100").Value = 16.75730000000001
$ws.Range("E102").Value = 16.67549999999999
